# Weekly update: insert a new record at row 22 (week of 2022-01-10) for
# "Terminal Hortofrutícola Agro Chillán - Arveja Verde", pushing all the
# older records that used to occupy rows 22-50 down by one row (to 23-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22; everything below (old rows 22-50)
# shifts down to rows 23-51, carrying its data and formatting with it.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with this week's new data.
$ws.Cells.Item(22, 1).Value = 7
$ws.Cells.Item(22, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(22, 3).Value = "Ñuble"
$ws.Cells.Item(22, 4).Value = 44571
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 16
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 27000
$ws.Cells.Item(22, 12).Value = 28000
$ws.Cells.Item(22, 13).Value = 27500
$ws.Cells.Item(22, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(22, 16).Value = 1100
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
